# Add list courses, statistical sourses of registeration
#
# The student sheet gained several new columns describing each student
# (phone number, day of birth, home town, gender, id number, current
# address). Update the header row so column C becomes "PhoneNumber",
# insert the new headers after "Email", and push "RegisterType" /
# "IdCourse" to the right (columns J and K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C used to be "Sdth" -> rename/replace with "PhoneNumber"
$ws.Range("C1").Value = "PhoneNumber"

# Columns E..I used to hold "Ngay sinh", "Que", "Gioi tinh", "CMND" and
# "Cho o hien tai" - replace them with their new English equivalents.
$ws.Range("E1").Value = "Day Of Birth"
$ws.Range("F1").Value = "HomeTown"
$ws.Range("G1").Value = "Gender"
$ws.Range("H1").Value = "IdNumber"
$ws.Range("I1").Value = "CurrentAddress"

# Put the focus/selection on K1, matching the saved view state.
$ws.Range("K1").Select()
